# Adds data for 2022-08-27 across the Citywide Totals, By Neighborhood,
# and per-neighborhood sheets in the CTA violent-crime YTD workbook.
$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = 47
$ws.Range("H2").Value = 72
$ws.Range("I2").Value = 79
$ws.Range("E3").Value = 97
$ws.Range("I3").Value = 136
$ws.Range("B6").Value = 256
$ws.Range("C6").Value = 328
$ws.Range("D6").Value = 296
$ws.Range("E6").Value = 287
$ws.Range("G6").Value = 337
$ws.Range("B7").Value = 350
$ws.Range("C7").Value = 442
$ws.Range("D7").Value = 460
$ws.Range("E7").Value = 438
$ws.Range("G7").Value = 487
$ws.Range("H7").Value = 470
$ws.Range("I7").Value = 600

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 8

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 7

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("C6").Value = 34
$ws.Range("C7").Value = 38

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("C2").Value = 1
$ws.Range("C5").Value = 11

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("E8").Value = 33
$ws.Range("I23").Value = 4
$ws.Range("C28").Value = 38
$ws.Range("B50").Value = 8
$ws.Range("B51").Value = 4
$ws.Range("H53").Value = 55
$ws.Range("I53").Value = 95
$ws.Range("I65").Value = 19
$ws.Range("D72").Value = 4
$ws.Range("G74").Value = 12
$ws.Range("G78").Value = 13
$ws.Range("C82").Value = 11
$ws.Range("G88").Value = 7
$ws.Range("E94").Value = 5
$ws.Range("C95").Value = 2
$ws.Range("B98").Value = 350
$ws.Range("C98").Value = 442
$ws.Range("D98").Value = 460
$ws.Range("E98").Value = 438
$ws.Range("G98").Value = 487
$ws.Range("H98").Value = 470
$ws.Range("I98").Value = 600

# Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("G4").Value = 12
$ws.Range("G5").Value = 13

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H2").Value = 7
$ws.Range("I3").Value = 23
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 95

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 4
$ws.Range("I6").Value = 19

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("G5").Value = 10
$ws.Range("G6").Value = 12

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 4

# Printers Row
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 4

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("B2").Value = 1
$ws.Range("B5").Value = 4

# West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 5

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E3").Value = 6
$ws.Range("E6").Value = 33

# Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2
